# "updates to single sat generator"
#
# 1. The shared "include_billboard=True))" string used by column J of the
#    obs_generator_billboard sheet gains a ",name_with_num=True" suffix.
# 2. Rows 15-36 of that sheet no longer have a populated "K" (generated
#    python line) column - their K cells are cleared out entirely.
# 3. Sheet selections / active tab change: equat_targets_1 stops being the
#    selected tab (and gets zoomed to 145%), obs_generator_billboard becomes
#    the selected tab with K2:K14 selected.

$wb = $excel.ActiveWorkbook

$equat = $wb.Worksheets.Item("equat_targets_1")
$billboard = $wb.Worksheets.Item("obs_generator_billboard")

# --- obs_generator_billboard: update the generated-code text -------------
# All of J2:J36 share the same underlying string; rewriting the whole
# original range in one shot keeps them sharing a single updated string
# (matching how the source file only edits the <t> of the existing shared
# string) instead of splintering off a second copy of the text.
$billboard.Range("J2:J36").Value = ",include_billboard=True,name_with_num=True))"

# Targets past row 14 (rows 15-36) no longer have a generated K value -
# clear those cells out entirely so they drop out of the sheet.
$billboard.Range("K15:K36").ClearContents()

# --- equat_targets_1: view changes ---------------------------------------
$equat.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 145
$equat.Range("B16").Select() | Out-Null

# --- obs_generator_billboard becomes the active/selected tab -------------
$billboard.Activate() | Out-Null
$billboard.Range("K2:K14").Select() | Out-Null
